$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new data row (17) mirroring the existing rows above it.
$srcRow = 16
$row = 17

# Copy formatting (including number format / style) from the row above
# so that we reuse the existing style definitions instead of creating new ones.
$ws.Range("A$srcRow`:N$srcRow").Copy()
$ws.Range("A$row`:N$row").PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item($row, 1).Value = 42622.888506944444
$ws.Cells.Item($row, 2).Value = 29
$ws.Cells.Item($row, 3).Value = 0
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = "Random"

$wb.Save()
